# Implements a purchase subsidy/ITC for clean industrial heat equipment
#
# Adds a new acronym-key entry (IFTQfS = "Industrial Fuels Qualifying for
# Subsidies") on the "Key to Variables" sheet, inserted right before the
# existing "ItUBB" row so the new row becomes row 191 and every row that
# used to start at 191 shifts down by one (through row 301).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Key to Variables")

# Insert a fresh blank row at 191; everything below (old 191..300) shifts
# down to 192..301.
$ws.Rows.Item(191).Insert()

# Copy the formatting (fill/alignment) of the "Importance" column from the
# row directly below (which holds the same "low" rating / style) so the new
# row's F cell reuses the existing style instead of inheriting the row
# above's style.
$ws.Cells.Item(192, 6).Copy()
$ws.Cells.Item(191, 6).PasteSpecial(-4122)

# Populate the new row's contents.
$ws.Cells.Item(191, 1).Value = "indst"
$ws.Cells.Item(191, 2).Value = "IFTQfS"
$ws.Cells.Item(191, 3).Value = "Industrial Fuels Qualifying for Subsidies"
$ws.Cells.Item(191, 4).Value = "Industrial Fuels Qualifying for PTC, Industrial Fuels Qualifying for ITC"
$ws.Cells.Item(191, 6).Value = "low"

# The inherited insert left a stray formatted-but-empty G cell; drop it so
# the row only carries A-D and F, matching the rest of the table.
$ws.Cells.Item(191, 7).Clear()

# Reflect the cell the author ended up on after making the edit.
$ws.Range("C191").Select()
